$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = 7.5
$ws.Range("U4").Value = 2.38
$ws.Range("V4").Value = 1.53
$ws.Range("Z4").Value = 9.5
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 101
$ws.Range("AZ4").Value = 201
$ws.Range("BA4").Value = 251

$ws.Range("G9").Value = 3.4
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 34
$ws.Range("AL9").Value = 29
$ws.Range("AR9").Value = 81
$ws.Range("AT9").Value = 2.63
$ws.Range("AW9").Value = 4.33

$ws.Range("G15").Value = 3.9
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 4.75
$ws.Range("S15").Value = 1.57
$ws.Range("T15").Value = 2.25
$ws.Range("U15").Value = 2.1
$ws.Range("V15").Value = 1.67
$ws.Range("W15").Value = 9
$ws.Range("X15").Value = 19
$ws.Range("Y15").Value = 15
$ws.Range("AA15").Value = 41
$ws.Range("AB15").Value = 51
$ws.Range("AC15").Value = 6.5
$ws.Range("AG15").Value = 5.5
$ws.Range("AH15").Value = 8.5
$ws.Range("AJ15").Value = 17
$ws.Range("AL15").Value = 41
$ws.Range("AP15").Value = 41
$ws.Range("AS15").Value = 401
$ws.Range("AT15").Value = 2.25
$ws.Range("AU15").Value = 9.5
$ws.Range("AV15").Value = 81
$ws.Range("AW15").Value = 3.75
$ws.Range("AY15").Value = 29
$ws.Range("BD15").Value = 151

$ws.Range("G16").Value = 2.75
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 2.38
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 3.1
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("Q16").Value = 1.98
$ws.Range("R16").Value = 1.88
$ws.Range("W16").Value = 9.5
$ws.Range("AC16").Value = 10
$ws.Range("AI16").Value = 10
$ws.Range("AL16").Value = 29
$ws.Range("AN16").Value = 4.75
$ws.Range("AO16").Value = 15
$ws.Range("AU16").Value = 8
$ws.Range("AY16").Value = 23
$ws.Range("BA16").Value = 67

$ws.Range("G17").Value = 1.95
$ws.Range("I17").Value = 3.6
$ws.Range("J17").Value = 2.75
$ws.Range("Q17").Value = 2.08
$ws.Range("R17").Value = 1.73
$ws.Range("AG17").Value = 10
$ws.Range("AK17").Value = 29
$ws.Range("AM17").Value = 251
$ws.Range("AU17").Value = 8
$ws.Range("BB17").Value = 201
